# Apply updated simulation results (pl_mw) for Case_3_250 (380 kV case)
# to worksheet "Sheet1" in the active workbook.
# Columns B..O map to the 14 data series, with E, I, M, N always 0.
# Only columns B, C, D, F, G, H, J, K, L, O change value; rows 2..25 (24 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","F","G","H","J","K","L","O")
$data = @(
    @(0.647781071271595,0.2467209802383792,0.2069910366683843,1.524588704530693,0.8672303508589039,0.9434404649970816,0.2369217653543636,0.3152905831755959,0.314115743149074,3.650875810708698),
    @(0.6045111744801943,0.2477725838488105,0.2043544428690467,1.53342034805209,0.8741143931177575,0.9505207086488667,0.2384331266282729,0.2784063110692045,0.3083792704219377,3.680035733597322),
    @(0.5780667066601382,0.2484695759305922,0.2028091948180588,1.539547951579763,0.8788168398963236,0.9552189054591196,0.2394481850584995,0.2557179949669575,0.3049705646521659,3.699673993939967),
    @(0.5673222571210772,0.2487665419038478,0.2021981087480285,1.542222435157441,0.8808527443470595,0.9572217881590674,0.2398837555892932,0.2464625982988338,0.3036101967021665,3.708112994470625),
    @(0.5655401003307361,0.2488166352291792,0.2020977656136225,1.542677252855185,0.8811980300644677,0.9575597029188714,0.2399574068979042,0.2449251783947233,0.3033860473154348,3.709540639390639),
    @(0.5779216730240933,0.2484735284932853,0.2028008779926154,1.539583301942045,0.8788438124538374,0.9552455592702316,0.2394539704976992,0.2555932119680762,0.3049521018011063,3.699786038582289),
    @(0.6328364895445304,0.2470729463743595,0.206066715681672,1.527487647918115,0.86950526691286,0.9458089600167128,0.2374248279031583,0.302581815806974,0.3121143439375516,3.660570411141052),
    @(0.7414718124299498,0.2447319221622593,0.2130513374722085,1.509354453083311,0.8549655583007763,0.9300839223771504,0.2341353049076282,0.3943754348819652,0.3270538930680118,3.597417061160925),
    @(0.8218293050953775,0.2432570943669319,0.2185321790585704,1.499428974112149,0.8465829913451302,0.9202199838543805,0.2321371927716065,0.4615768292617588,0.3385682376139556,3.559387409143483),
    @(0.8584971369661787,0.242638946453539,0.2211004704087713,1.495649478131575,0.8432688525889631,0.9160982558374968,0.2313187587549628,0.4920915202359311,0.3439218897598124,3.54390150787097),
    @(0.8723977601972877,0.2424124221233228,0.2220837126182573,1.494323912425678,0.8420856408522752,0.9145899234722776,0.23102182657966,0.5036380870089658,0.3459656713590675,3.538298049654401),
    @(0.8694033457773003,0.2424608726877864,0.2218714801940251,1.494604700237304,0.842337274117483,0.9149124370576089,0.2310851988583948,0.5011517245447976,0.3455247766436003,3.539493261303932),
    @(0.8596404479547459,0.2426201589302508,0.2211811487636339,1.495538306508685,0.8431700704234402,0.9159731131248066,0.2312940697508346,0.4930416412556724,0.344089703800762,3.543435283598967),
    @(0.8536623581187825,0.2427187092825207,0.2207596894060799,1.496123921589422,0.8436895303375707,0.9166296397286544,0.2314237002704544,0.4880728295084111,0.3432128200311979,3.545883835152836),
    @(0.8194351681294165,0.2432985507018657,0.2183658361432776,1.499690764454172,0.8468096224642352,0.9204966956904599,0.2321924986377226,0.4595814449714339,0.3382206778398142,3.560435939878602),
    @(0.7984661396054378,0.2436677546755845,0.2169164280901157,1.502067227120094,0.848851542314712,0.9229625530106418,0.2326872971878906,0.4420881779009562,0.3351876773707829,3.569827694890876),
    @(0.7864159810660851,0.2438850788064499,0.2160898351613127,1.503503357172008,0.850072981049756,0.924415248893979,0.2329804136414424,0.4320213251554605,0.3334540779404023,3.575400332076413),
    @(0.7823378651380608,0.243959515214506,0.2158111821168092,1.504001504378977,0.8504946079634408,0.9249130171692812,0.2330811221673912,0.4286119919194675,0.3328689886239715,3.577316461422257),
    @(0.8006972305708473,0.2436279383683981,0.2170699892622707,1.501807082337194,0.8486293141155201,0.9226964985395583,0.2326337432863106,0.4439509080445418,0.335509418346092,3.568810255696349),
    @(0.8625076394067435,0.2425731679143155,0.2213836266445384,1.495261217584542,0.8429235098098076,0.9156601432204496,0.231232366935096,0.4954240105298595,0.3445107738568538,3.542270341374518),
    @(0.9029931385289274,0.2419278355222261,0.2242650765027321,1.491598845678958,0.8396128319097755,0.9113673190542642,0.2303921964379434,0.5290137180011243,0.3504895990608503,3.526444566986953),
    @(0.8813774385438933,0.2422682440354222,0.2227215317060001,1.493497229220601,0.8413415180380639,0.9136305191493648,0.23083369210255,0.5110911490860985,0.3472898672665394,3.534752081273609),
    @(0.7996885380324272,0.2436459235459694,0.2170005434070674,1.501924476143898,0.8487296354891143,0.9228166726078939,0.2326579280588525,0.4431087980352686,0.3353639277385412,3.569269700465185),
    @(0.7119852163426117,0.2453220356026407,0.2111001541950088,1.513662794567047,0.8584950383579297,0.9340409107028975,0.2349515565776237,0.3695830353864551,0.3229173201078055,3.613031124337184)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $cellRef = "$($cols[$j])$rowNum"
        $ws.Range($cellRef).Value = $rowValues[$j]
    }
}
